# Planification Initiale DMA.pdf ajoute / Fin de la planification initiale
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekday header strings (row 1, columns D:H) -> appended to sharedStrings
$ws.Range("D1").Value = "lu"
$ws.Range("E1").Value = "ma"
$ws.Range("F1").Value = "me"
$ws.Range("G1").Value = "je"
$ws.Range("H1").Value = "ve"

# Row 2 (week 1) - updated daily hours
$ws.Range("E2").Value = 0.2986111111111111   # 7:10
$ws.Range("F2").Value = 0.16666666666666666  # 4:00
$ws.Range("G2").Value = 0.2986111111111111   # 7:10

# Row 3 (week 2) - updated daily hours
$ws.Range("D3").Value = 0.23263888888888887  # 5:35
$ws.Range("E3").Value = 0.2986111111111111   # 7:10
$ws.Range("F3").Value = 0.16666666666666666  # 4:00
$ws.Range("G3").Value = 0.2986111111111111   # 7:10

# Row 4 (week 3) - updated daily hours
$ws.Range("D4").Value = 0.23263888888888887  # 5:35
$ws.Range("E4").Value = 0.2986111111111111   # 7:10
$ws.Range("F4").Value = 0.16666666666666666  # 4:00

# Row 5 (week 4) - updated daily hours
$ws.Range("D5").Value = 0.23263888888888887  # 5:35
$ws.Range("E5").Value = 0.2986111111111111   # 7:10
$ws.Range("F5").Value = 0.16666666666666666  # 4:00
$ws.Range("G5").Value = 0.2986111111111111   # 7:10

# Row 6 (week 5) - updated daily hours + total now references only ma
$ws.Range("E6").Value = 0.2986111111111111   # 7:10
$ws.Range("F6").Value = 0.16666666666666666  # 4:00
$ws.Range("I6").Formula = "=E6"

# Row 7 (week 6) - updated daily hours + total forced to 0
$ws.Range("D7").Value = 0.23263888888888887  # 5:35
$ws.Range("E7").Value = 0.2986111111111111   # 7:10
$ws.Range("F7").Value = 0.16666666666666666  # 4:00
$ws.Range("I7").Formula = "=0"

# Row 8 (week 7) - grand total entered directly (plain 2-decimal number, not h:mm)
$ws.Range("I8").NumberFormat = "0.00"
$ws.Range("I8").Value = 90.05

# Printing setup for the exported "Planification Initiale DMA.pdf"
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final cursor position left on the sheet
$ws.Range("F28").Select() | Out-Null
